$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Anxa2"
$ws.Range("C2").Value = "Tlr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 93.08856299999998
$ws.Range("H2").Value = 279.265689
$ws.Range("I2").Value = 0.1947836339852847
$ws.Range("J2").Value = 0.1947836339852847
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 44.04866033333334
$ws.Range("N2").Value = 132.145981
$ws.Range("O2").Value = 0.3636142564479216
$ws.Range("P2").Value = 0.3636142564479216
$ws.Range("Q2").Value = 4100.4264925051
$ws.Range("R2").Value = 36903.8384325459
$ws.Range("S2").Value = 0.0708261062397834
$ws.Range("T2").Value = 0.0708261062397834

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Anxa2"
$ws.Range("C3").Value = "Tlr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 93.08856299999998
$ws.Range("H3").Value = 279.265689
$ws.Range("I3").Value = 0.1947836339852847
$ws.Range("J3").Value = 0.1947836339852847
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.226320666666667
$ws.Range("N3").Value = 6.678962
$ws.Range("O3").Value = 0.01837790134135009
$ws.Range("P3").Value = 0.01837790134135009
$ws.Range("Q3").Value = 207.244991637202
$ws.Range("R3").Value = 1865.204924734818
$ws.Range("S3").Value = 0.003579714408291208
$ws.Range("T3").Value = 0.003579714408291208

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Anxa2"
$ws.Range("C4").Value = "Tlr2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 93.08856299999998
$ws.Range("H4").Value = 279.265689
$ws.Range("I4").Value = 0.1947836339852847
$ws.Range("J4").Value = 0.1947836339852847
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 74.86619966666667
$ws.Range("N4").Value = 224.598599
$ws.Range("O4").Value = 0.6180078422107282
$ws.Range("P4").Value = 0.6180078422107284
$ws.Range("Q4").Value = 6969.186944241078
$ws.Range("R4").Value = 62722.68249816971
$ws.Range("S4").Value = 0.1203778133372101
$ws.Range("T4").Value = 0.1203778133372101

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Anxa2"
$ws.Range("C5").Value = "Tlr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 153.2725883333333
$ws.Range("H5").Value = 459.817765
$ws.Range("I5").Value = 0.320716001877666
$ws.Range("J5").Value = 0.3207160018776659
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 44.04866033333334
$ws.Range("N5").Value = 132.145981
$ws.Range("O5").Value = 0.3636142564479216
$ws.Range("P5").Value = 0.3636142564479216
$ws.Range("Q5").Value = 6751.452181905831
$ws.Range("R5").Value = 60763.06963715247
$ws.Range("S5").Value = 0.1166169105536977
$ws.Range("T5").Value = 0.1166169105536977

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Anxa2"
$ws.Range("C6").Value = "Tlr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 153.2725883333333
$ws.Range("H6").Value = 459.817765
$ws.Range("I6").Value = 0.320716001877666
$ws.Range("J6").Value = 0.3207160018776659
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.226320666666667
$ws.Range("N6").Value = 6.678962
$ws.Range("O6").Value = 0.01837790134135009
$ws.Range("P6").Value = 0.01837790134135009
$ws.Range("Q6").Value = 341.2339310399922
$ws.Range("R6").Value = 3071.10537935993
$ws.Range("S6").Value = 0.005894087041099994
$ws.Range("T6").Value = 0.005894087041099993

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Anxa2"
$ws.Range("C7").Value = "Tlr2"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 153.2725883333333
$ws.Range("H7").Value = 459.817765
$ws.Range("I7").Value = 0.320716001877666
$ws.Range("J7").Value = 0.3207160018776659
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 74.86619966666667
$ws.Range("N7").Value = 224.598599
$ws.Range("O7").Value = 0.6180078422107282
$ws.Range("P7").Value = 0.6180078422107284
$ws.Range("Q7").Value = 11474.93620159014
$ws.Range("R7").Value = 103274.4258143113
$ws.Range("S7").Value = 0.1982050042828682
$ws.Range("T7").Value = 0.1982050042828682

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Anxa2"
$ws.Range("C8").Value = "Tlr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 97.56176266666667
$ws.Range("H8").Value = 292.685288
$ws.Range("I8").Value = 0.2041436032289296
$ws.Range("J8").Value = 0.2041436032289296
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 44.04866033333334
$ws.Range("N8").Value = 132.145981
$ws.Range("O8").Value = 0.3636142564479216
$ws.Range("P8").Value = 0.3636142564479216
$ws.Range("Q8").Value = 4297.464945225282
$ws.Range("R8").Value = 38677.18450702753
$ws.Range("S8").Value = 0.07422952449668675
$ws.Range("T8").Value = 0.07422952449668675

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Anxa2"
$ws.Range("C9").Value = "Tlr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 97.56176266666667
$ws.Range("H9").Value = 292.685288
$ws.Range("I9").Value = 0.2041436032289296
$ws.Range("J9").Value = 0.2041436032289296
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.226320666666667
$ws.Range("N9").Value = 6.678962
$ws.Range("O9").Value = 0.01837790134135009
$ws.Range("P9").Value = 0.01837790134135009
$ws.Range("Q9").Value = 217.2037685012284
$ws.Range("R9").Value = 1954.833916511056
$ws.Range("S9").Value = 0.003751730999608985
$ws.Range("T9").Value = 0.003751730999608985

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Anxa2"
$ws.Range("C10").Value = "Tlr2"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 97.56176266666667
$ws.Range("H10").Value = 292.685288
$ws.Range("I10").Value = 0.2041436032289296
$ws.Range("J10").Value = 0.2041436032289296
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 74.86619966666667
$ws.Range("N10").Value = 224.598599
$ws.Range("O10").Value = 0.6180078422107282
$ws.Range("P10").Value = 0.6180078422107284
$ws.Range("Q10").Value = 7304.078403634613
$ws.Range("R10").Value = 65736.70563271153
$ws.Range("S10").Value = 0.1261623477326339
$ws.Range("T10").Value = 0.1261623477326339

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Anxa2"
$ws.Range("C11").Value = "Tlr2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 133.9846036666667
$ws.Range("H11").Value = 401.953811
$ws.Range("I11").Value = 0.2803567609081197
$ws.Range("J11").Value = 0.2803567609081197
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 44.04866033333334
$ws.Range("N11").Value = 132.145981
$ws.Range("O11").Value = 0.3636142564479216
$ws.Range("P11").Value = 0.3636142564479216
$ws.Range("Q11").Value = 5901.842296809288
$ws.Range("R11").Value = 53116.58067128359
$ws.Range("S11").Value = 0.1019417151577537
$ws.Range("T11").Value = 0.1019417151577537

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Anxa2"
$ws.Range("C12").Value = "Tlr2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 133.9846036666667
$ws.Range("H12").Value = 401.953811
$ws.Range("I12").Value = 0.2803567609081197
$ws.Range("J12").Value = 0.2803567609081197
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.226320666666667
$ws.Range("N12").Value = 6.678962
$ws.Range("O12").Value = 0.01837790134135009
$ws.Range("P12").Value = 0.01837790134135009
$ws.Range("Q12").Value = 298.2926921582424
$ws.Range("R12").Value = 2684.634229424182
$ws.Range("S12").Value = 0.005152368892349899
$ws.Range("T12").Value = 0.005152368892349898

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Anxa2"
$ws.Range("C13").Value = "Tlr2"
$ws.Range("D13").Value = "M2"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 133.9846036666667
$ws.Range("H13").Value = 401.953811
$ws.Range("I13").Value = 0.2803567609081197
$ws.Range("J13").Value = 0.2803567609081197
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 74.86619966666667
$ws.Range("N13").Value = 224.598599
$ws.Range("O13").Value = 0.6180078422107282
$ws.Range("P13").Value = 0.6180078422107284
$ws.Range("Q13").Value = 10030.91809036787
$ws.Range("R13").Value = 90278.2628133108
$ws.Range("S13").Value = 0.1732626768580161
$ws.Range("T13").Value = 0.1732626768580161

